$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's data as row 94 (daily auto push)
$row = 94

# Use a leading apostrophe so the date-like text "2025/10/12" is stored
# as literal text (matching the source data format) instead of being
# auto-converted to a date serial number; ClearFormats() strips the
# quote-prefix cell style iron_native applies so the cell keeps the
# workbook's default (unstyled) formatting, same as the rest of the row.
$ws.Cells.Item($row, 1).Value = "'2025/10/12"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "日"
$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 201
